$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G (VAD-Wiener-GSC / Fwiener-GSC results) gets an explicit width.
$ws.Columns.Item(7).ColumnWidth = 13.4

# Row 1 merged-banner row grows two blank, styled cells into F1/G1.
$ws.Range("F1").NumberFormat = "@"
$ws.Range("G1").NumberFormat = "@"

# Make sure every new cell we touch is stored as text (matches the existing
# "@" text format used throughout columns B:E), so numeric-looking labels
# like "-7.72" aren't coerced into real numbers.
$ws.Range("F2:G11").NumberFormat = "@"

# --- Column G: Fwiener-GSC -------------------------------------------------
$ws.Range("G2").Value = "Fwiener-GSC"
$ws.Range("G9").Value = "-7.72"
$ws.Range("G11").Value = "5.70"
$ws.Range("G10").Value = "-0.19"
$ws.Range("G8").Value = "5.97"
$ws.Range("G7").Value = "0.89"
$ws.Range("G6").Value = "-5.81"
$ws.Range("G3").Value = "-3.14"
$ws.Range("G4").Value = "3.73"
$ws.Range("G5").Value = "9.01"

# --- Column F: VAD-Wiener-GSC -----------------------------------------------
$ws.Range("F5").Value = "12.00"
$ws.Range("F4").Value = "7.15"
$ws.Range("F3").Value = "2.37"
$ws.Range("F2").Value = "VAD-Wiener-GSC"
$ws.Range("F6").Value = "0.01"
$ws.Range("F7").Value = "2.22"
$ws.Range("F8").Value = "7.28"
$ws.Range("F11").Value = "9.48"
$ws.Range("F10").Value = "4.24"
$ws.Range("F9").Value = "0.73"

# Author's cursor ended on F9 after entering the last value.
$ws.Range("F9").Select()
